$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J10").Value = 3
$ws.Range("J11").Value = 2
$ws.Range("J14").Value = "Roblez"
$ws.Range("J16").Value = "88-90 MPH"
$ws.Range("J17").Value = "CB,FB,CH"
$ws.Range("J19").Value = 7
$ws.Range("M19").Value = "101.29 MPH"
$ws.Range("M21").Value = "11.83°"
$ws.Range("J23").Value = "Plum"
$ws.Range("M23").Value = "Line Drive"
$ws.Range("M24").Value = "Single"
$ws.Range("J25").Value = "84-86 MPH"
$ws.Range("J26").Value = "SL,FB,CH"
$ws.Range("J28").Value = 4
$ws.Range("M28").Value = "nan MPH"
$ws.Range("M30").Value = "nan°"
$ws.Range("J32").Value = "Herbst"
$ws.Range("M32").Value = "Undefined"
$ws.Range("J33").Value = "Right"
$ws.Range("M33").Value = "Undefined"
$ws.Range("J34").Value = "83-85 MPH"
$ws.Range("J35").Value = "SL,CB,FB,CH"
$ws.Range("J37").Value = 8
$ws.Range("M37").Value = "55.63 MPH"
$ws.Range("J38").Value = 1
$ws.Range("M39").Value = "-43.57°"
$ws.Range("J41").Value = "Thompson"
$ws.Range("M41").Value = "Ground Ball"
$ws.Range("J42").Value = "Left"
$ws.Range("M42").Value = "Out"
$ws.Range("J43").Value = "84-84 MPH"
$ws.Range("J44").Value = "SL,FB,CH"
$ws.Range("J46").Value = 1
$ws.Range("M46").Value = "93.8 MPH"
$ws.Range("J47").Value = 1
$ws.Range("M48").Value = "39.82°"
$ws.Range("M50").Value = "Fly Ball"
$ws.Range("M51").Value = "Out"
$ws.Range("J52").Value = "88-90 MPH"
$ws.Range("J53").Value = "CB,FB,CH"
$ws.Range("J61").Value = 6
$ws.Range("M61").Value = "nan MPH"
$ws.Range("J62").Value = 0
$ws.Range("M63").Value = "nan°"
$ws.Range("J65").Value = "Herbst"
$ws.Range("M65").Value = "Undefined"
$ws.Range("M66").Value = "Undefined"
$ws.Range("J67").Value = "83-85 MPH"
$ws.Range("J68").Value = "SL,CB,FB,CH"
